$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new log entry for row 6 (matches style of existing entries in rows 2 and 4)
$ws.Range("B6").Value = "Eetu"
$ws.Range("C6").Value = 1.7
$ws.Range("D6").Value = "JSON parsing. Dynamic views for Meetings. V1.1 release."

# Move the active selection to D7, as it was left after entering the new row
$ws.Range("D7").Select()
